$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) NewTreeView row (row 96) - mark columns C/D as implemented ("Y")
$ws.Cells.Item(96, 3).Value = "Y"
$ws.Cells.Item(96, 4).Value = "Y"

# 2) Insert a new row right after WriteTableToKml (row 221) for the new
#    WriteTableToShapefile command, pushing everything below down by one row.
$ws.Rows.Item(222).Insert()
$ws.Cells.Item(222, 1).Value = "WriteTableToShapefile"

# 3) Update the description of WriteTableToKml (row 221, column B)
$ws.Cells.Item(221, 2).Value = "Write a table to KML file, formatting spatial data columns into KML."

$ws.Cells.Item(222, 2).Value = "Write a table to Esri shapefile."
$ws.Cells.Item(222, 3).Value = "Y"
$ws.Cells.Item(222, 4).Value = "Y"

# Match the centered style ("s=3": horizontal-center) used throughout
# column C/D for the new cells (xlCenter = -4108).
$ws.Cells.Item(222, 3).HorizontalAlignment = -4108
$ws.Cells.Item(222, 4).HorizontalAlignment = -4108

# Refresh the view state to match where the edit left the selection.
$ws.Range("D223").Select()
